$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 with Problem Name and Platform entries
$ws.Range("C8").Value = "Maximum Gap"
$ws.Range("D8").Value = "Bosscoder Academy"

# Update the active selection to C9, matching the saved workbook view state
$ws.Range("C9").Select()
